# Daily attendance processing - 2026-01-09 08:42:58
# Swap the order of the "Recorded By" names in column G from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Only the "Recorded By" column (G) holds this value; walk the matches with
# Find/FindNext so only the affected cells are touched.
$col = $ws.Columns.Item(7)
$first = $col.Find($oldValue, [System.Reflection.Missing]::Value, -4163, 1, 1, 1, $true)

if ($first -ne $null) {
    $firstAddress = $first.Address()
    $cell = $first
    do {
        $cell.Value2 = $newValue
        $cell = $col.FindNext($cell)
    } while (($cell -ne $null) -and ($cell.Address() -ne $firstAddress))
}
